# Updated cryptos list on Thu Dec 28 10:49:45 UTC 2023 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns with new scrape values,
# and two rows (46/47 - BitcoinSV / Celestia) swapped rank position.
#
# Note: several Price values look numeric (e.g. "1.60", "9.34") but must be
# stored as literal text (to match the source feed's formatting, trailing
# zeros included). A plain Range.Value assignment would coerce those to
# numbers, so such values are written with a leading "'" (text/quote
# prefix) and the cell style is then reset to "Normal" so no quotePrefix
# formatting flag is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.278.66"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "2.404.39"
$ws.Range("E3").Value = "  +5.60%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").Value = "'328.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.35%  "
$ws.Range("D6").Value = "'106.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.70%  "
$ws.Range("D7").Value = "'0.654"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.23%  "
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +6.23%  "
$ws.Range("D10").Value = "'42.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.65%  "
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("E12").Value = "  -1.71%  "
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'17.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.23%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "2.764.88"
$ws.Range("E16").Value = "  +5.68%  "
$ws.Range("D17").Value = "2.405.77"
$ws.Range("E17").Value = "  +5.69%  "
$ws.Range("D18").Value = "43.199.74"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("E19").Value = "  +7.36%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "'77.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "'3.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.00%  "
$ws.Range("D23").Value = "'274.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.63%  "
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("D25").Value = "'9.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.07%  "
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'23.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("D29").Value = "'176.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.61%  "
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("D31").Value = "'37.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("D32").Value = "'0.0943"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.12%  "
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "'5.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.26%  "
$ws.Range("D35").Value = "'0.136"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.48%  "
$ws.Range("D36").Value = "'4.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("E40").Value = "  +15.62%  "
$ws.Range("D41").Value = "'1.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.08%  "
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").Value = "'70.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "'123.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +14.26%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'12.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.44%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "'91.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +41.35%  "
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").Value = "'9.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.53%  "
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").Value = "'0.496"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.73%  "
